$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 47, shifting rows 47:137 down to 48:138
$ws.Rows("47").Insert()

# Populate the newly inserted row 47 with the new weekly record
$ws.Range("A47").Value = 8
$ws.Range("B47").Value = "Terminal La Palmera de La Serena"
$ws.Range("C47").Value = "Coquimbo"
$ws.Range("D47").Value = 44533
$ws.Range("E47").Value = 4
$ws.Range("F47").Value = 100112037
$ws.Range("G47").Value = "Cebollín"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 3000
$ws.Range("K47").Value = 900
$ws.Range("L47").Value = 1000
$ws.Range("M47").Value = 950
$ws.Range("N47").Value = '$/paquete 6 unidades'
$ws.Range("O47").Value = "Provincia del Elquí"
$ws.Range("P47").Value = 158
$ws.Range("Q47").Value = 6
$ws.Range("R47").Value = "Hortaliza"

# Ensure the date cell keeps the date number format used by the other rows
$ws.Range("D47").NumberFormat = $ws.Range("D48").NumberFormat
